$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'
$ws.Range("B10").Value = 'Bejucal De Ocampo'
$ws.Range("B14").Value = 'Comitán De Domínguez'
$ws.Range("B24").Value = 'Mazapa De Madero'
$ws.Range("D24").Value = 0.009633911368015412
$ws.Range("B28").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B31").Value = 'San Cristóbal De Las Casas'
$ws.Range("B40").Value = 'Hidalgo Del Parral'
$ws.Range("A45").Value = 'Ciudad De México'
$ws.Range("D50").Value = 0.009633911368015412
$ws.Range("A62").Value = 'Estado De México'
$ws.Range("B65").Value = 'Ecatepec De Morelos'
$ws.Range("B67").Value = 'Ixtapan De La Sal'
$ws.Range("B70").Value = 'San Felipe Del Progreso'
$ws.Range("D73").Value = 0.009633911368015412
$ws.Range("B75").Value = 'Apaseo El Grande'
$ws.Range("B78").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B83").Value = 'San Francisco Del Rincón'
$ws.Range("B85").Value = 'San Luis De La Paz'
$ws.Range("B86").Value = 'Silao De La Victoria'
$ws.Range("B90").Value = 'Acapulco De Juárez'
$ws.Range("B93").Value = 'Ayutla De Los Libres'
$ws.Range("B94").Value = 'Chilapa De Álvarez'
$ws.Range("B96").Value = 'Coyuca De Catalán'
$ws.Range("B99").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B100").Value = 'Iguala De La Independencia'
$ws.Range("B102").Value = 'Mártir De Cuilapan'
$ws.Range("B105").Value = 'Taxco De Alarcón'
$ws.Range("B108").Value = 'Tlapa De Comonfort'
$ws.Range("B117").Value = 'Nopala De Villagrán'
$ws.Range("B124").Value = 'Autlán De Navarro'
$ws.Range("B129").Value = 'Cuautitlán De García Barragán'
$ws.Range("B134").Value = 'Jilotlán De Los Dolores'
$ws.Range("B135").Value = 'La Manzanilla De La Paz'
$ws.Range("B136").Value = 'San Miguel El Alto'
$ws.Range("B137").Value = 'Tizapán El Alto'
$ws.Range("B138").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B141").Value = 'Unión De Tula'
$ws.Range("D141").Value = 0.009633911368015412
$ws.Range("D143").Value = 0.09248554913294796
$ws.Range("B165").Value = 'Puente De Ixtla'
$ws.Range("B166").Value = 'Tetela Del Volcán'
$ws.Range("D173").Value = 0.009633911368015412
$ws.Range("D174").Value = 0.009633911368015412
$ws.Range("B175").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B176").Value = 'Oaxaca De Juárez'
$ws.Range("D176").Value = 0.009633911368015412
$ws.Range("B177").Value = 'Pinotepa De Don Luis'
$ws.Range("B178").Value = 'Putla Villa De Guerrero'
$ws.Range("B182").Value = 'San Dionisio Del Mar'
$ws.Range("B186").Value = 'San Miguel Del Puerto'
$ws.Range("B188").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B197").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B198").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B199").Value = 'Tlacolula De Matamoros'
$ws.Range("B200").Value = 'Totontepec Villa De Morelos'
$ws.Range("B212").Value = 'Izúcar De Matamoros'
$ws.Range("B225").Value = 'Tuzamapan De Galeana'
$ws.Range("B232").Value = 'San Juan Del Río'
$ws.Range("B240").Value = 'Villa De Reyes'
$ws.Range("B268").Value = 'Cosamaloapan De Carpio'
$ws.Range("B271").Value = 'Ignacio De La Llave'
$ws.Range("B276").Value = 'Martínez De La Torre'
$ws.Range("B278").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B280").Value = 'Poza Rica De Hidalgo'
$ws.Range("B284").Value = 'Soledad De Doblado'
$ws.Range("D288").Value = 0.009633911368015412

# --- Remove trailing footnote rows 296-301 (dimension shrinks to A1:D295) ---
$ws.Rows("296:301").Delete()
